$d = $word.ActiveDocument

$d.Content.Find.Execute("Warm Up", $true, $false, $false, $false, $false, $true, 1, $false, "Warm Up    5-10 minutes", 2)
$d.Content.Find.Execute("Squat", $true, $false, $false, $false, $false, $true, 1, $false, "Squat       4 sets", 2)
$d.Content.Find.Execute("Lunge", $true, $false, $false, $false, $false, $true, 1, $false, "Lunge       4 sets", 2)
$d.Content.Find.Execute("Stretch", $true, $false, $false, $false, $false, $true, 1, $false, "Stretch      5 minutes", 2)
